$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in "Вариант" (L column) values for students that were missing them ---
# Row 7  = Выборнов Даниил
# Row 14 = Иванова Снежана
# Row 17 = Кудрявцева Полина
# Row 21 = Пушкина Софья
# Row 27 = Тикконен Герман
# Row 30 = Хабибулина Майя
# Row 33 = Юшина Полина

# L7 matches the formatting already used by plain (unstyled) "Вариант" cells like L5/L12
$ws.Range("L5").Copy($ws.Range("L7"))
$ws.Range("L7").Value = 4

# L14, L17, L21, L27, L30, L33 match the formatting used by L9/L16/L18/L23 (centered, wrapped)
$ws.Range("L9").Copy($ws.Range("L14"))
$ws.Range("L14").Value = 2

$ws.Range("L9").Copy($ws.Range("L17"))
$ws.Range("L17").Value = 4

$ws.Range("L9").Copy($ws.Range("L21"))
$ws.Range("L21").Value = 2

$ws.Range("L9").Copy($ws.Range("L27"))
$ws.Range("L27").Value = 4

$ws.Range("L9").Copy($ws.Range("L30"))
$ws.Range("L30").Value = 1

$ws.Range("L9").Copy($ws.Range("L33"))
$ws.Range("L33").Value = 2

# --- Add a new student row 35 (a late addition to the roster) ---
# Formatting matches the "Вариант" header cell (L3): bold, centered.
$ws.Range("L3").Copy($ws.Range("B35"))
$ws.Range("B35").Value = "Даниил романович"

# --- Update the active selection to reflect where the editor last worked ---
$null = $ws.Range("L17").Select()
